# Append: 2026-01-28 18:42 JST
# Update the "取得日時" (retrieved-at) timestamp in column A for the
# currently-listed rows (2-9) on the "ランサーズ" sheet to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-28 18:42:35"

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
